$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.0082
$ws.Range("I4").Value = -0.0365
$ws.Range("J4").Value = -0.0354
$ws.Range("N4").Value = -0.0073
$ws.Range("P4").Value = -0.0687
$ws.Range("Q4").Value = -0.0299
$ws.Range("I5").Value = -0.1777
$ws.Range("J5").Value = -0.4306
$ws.Range("I6").Value = -0.0587
$ws.Range("J6").Value = -0.0852
$ws.Range("K6").Value = -0.0915
$ws.Range("L6").Value = -0.076
$ws.Range("M6").Value = -0.0726
$ws.Range("N6").Value = -0.0625
$ws.Range("O6").Value = -0.0492
$ws.Range("P6").Value = -0.0529
$ws.Range("Q6").Value = -0.0338
$ws.Range("I7").Value = -0.0735
$ws.Range("J7").Value = -0.0839
$ws.Range("K7").Value = 0.0636
$ws.Range("L7").Value = -0.0358
$ws.Range("M7").Value = -0.0213
$ws.Range("N7").Value = -0.0458
$ws.Range("O7").Value = -0.0706
$ws.Range("P7").Value = 0.0349
$ws.Range("Q7").Value = 0.0334
$ws.Range("I8").Value = -0.453
$ws.Range("J8").Value = -0.2229
$ws.Range("K8").Value = -0.0924
$ws.Range("L8").Value = -0.0679
$ws.Range("M8").Value = 0.0287
$ws.Range("N8").Value = 0.3272
$ws.Range("O8").Value = 0.321
$ws.Range("P8").Value = 0.3382
$ws.Range("Q8").Value = 0.3453
$ws.Range("I9").Value = -0.0907
$ws.Range("J9").Value = -0.0624
$ws.Range("K9").Value = -0.0388
$ws.Range("L9").Value = -0.0157
$ws.Range("M9").Value = -0.0178
$ws.Range("N9").Value = -0.0394
$ws.Range("O9").Value = -0.0344
$ws.Range("P9").Value = -0.0284
$ws.Range("I10").Value = -0.0927
$ws.Range("J10").Value = -0.0806
$ws.Range("K10").Value = -0.0939
$ws.Range("L10").Value = -0.0927
$ws.Range("M10").Value = -0.0817
$ws.Range("I11").Value = 0.0071
$ws.Range("I12").Value = -0.5943
$ws.Range("J12").Value = -0.1432
$ws.Range("K12").Value = -0.0174
$ws.Range("L12").Value = 0.0456
$ws.Range("I14").Value = -0.0708
$ws.Range("J14").Value = -0.0589
$ws.Range("K14").Value = -0.5233
$ws.Range("L14").Value = -0.5717
$ws.Range("M14").Value = -0.2222
$ws.Range("N14").Value = -0.1758
$ws.Range("O14").Value = -0.3206
$ws.Range("P14").Value = -0.2399
$ws.Range("Q14").Value = -0.1032
$ws.Range("I15").Value = -0.3956
$ws.Range("J15").Value = -0.3902
$ws.Range("K15").Value = -0.2835
$ws.Range("L15").Value = -0.1285
$ws.Range("M15").Value = -0.0067
$ws.Range("N15").Value = -0.0009
$ws.Range("O15").Value = -0.0002
$ws.Range("P15").Value = 0.0001
$ws.Range("Q15").Value = 0.0001
$ws.Range("I16").Value = -2.1018
$ws.Range("J16").Value = -1.8326
$ws.Range("K16").Value = -2.5519
$ws.Range("L16").Value = -1.0643
$ws.Range("M16").Value = -0.4702
$ws.Range("N16").Value = -0.1682
$ws.Range("O16").Value = -0.368
$ws.Range("P16").Value = -0.3384
$ws.Range("Q16").Value = -0.0074
$ws.Range("I18").Value = -0.1803
$ws.Range("J18").Value = -0.4223
$ws.Range("I22").Value = 0.1361
$ws.Range("J22").Value = 0.3552
$ws.Range("I23").Value = -0.0064
$ws.Range("K23").Value = -0.0086
$ws.Range("L23").Value = -0.0001
$ws.Range("N23").Value = 0.0039
$ws.Range("P23").Value = 0.0064
$ws.Range("I24").Value = -0.0786
$ws.Range("J24").Value = 0.0049
$ws.Range("K24").Value = -0.0816
$ws.Range("L24").Value = 0.1253
$ws.Range("M24").Value = 0.111
$ws.Range("N24").Value = 0.0968
$ws.Range("O24").Value = 0.083
$ws.Range("P24").Value = -0.0598
$ws.Range("Q24").Value = -0.0478
$ws.Range("I25").Value = 0.004
$ws.Range("J25").Value = -0.0043
$ws.Range("K25").Value = 0.0616
$ws.Range("L25").Value = 0.0007
$ws.Range("M25").Value = 0.0194
$ws.Range("N25").Value = 0.0075
$ws.Range("O25").Value = -0.0079
$ws.Range("P25").Value = -0.0151
$ws.Range("J29").Value = -0.0403
$ws.Range("K29").Value = -0.0274
$ws.Range("L29").Value = -0.0163
$ws.Range("I30").Value = 0.0014
$ws.Range("I32").Value = -0.0001
$ws.Range("J32").Value = -0.0001
$ws.Range("P32").Value = 0
$ws.Range("I33").Value = 0.038
$ws.Range("J33").Value = -0.0728
$ws.Range("I34").Value = -0.0029
$ws.Range("J34").Value = -0.0029
$ws.Range("K34").Value = -0.0029
$ws.Range("L34").Value = -0.003
$ws.Range("M34").Value = -0.003
$ws.Range("N34").Value = -0.003
$ws.Range("O34").Value = -0.0029
$ws.Range("P34").Value = -0.0029
$ws.Range("Q34").Value = -0.0029
$ws.Range("I35").Value = -0.0181
$ws.Range("J35").Value = -0.0198
$ws.Range("K35").Value = -0.0216
$ws.Range("L35").Value = -0.0208
$ws.Range("M35").Value = -0.0038
$ws.Range("N35").Value = -0.0017
$ws.Range("O35").Value = 0.0004
$ws.Range("P35").Value = 0.0006
$ws.Range("Q35").Value = 0.0008
$ws.Range("I36").Value = 0.0563
$ws.Range("J36").Value = 0.0294
$ws.Range("K36").Value = 0.007
$ws.Range("L36").Value = 0.0168
$ws.Range("M36").Value = 0.0164
$ws.Range("N36").Value = 0.0159
$ws.Range("O36").Value = 0.0157
$ws.Range("P36").Value = 0.0105
$ws.Range("Q36").Value = -0.0104
$ws.Range("I37").Value = 0.0111
$ws.Range("J37").Value = -0.0032
$ws.Range("K37").Value = -0.0001
$ws.Range("L37").Value = -0.0039
$ws.Range("M37").Value = -0.0001
$ws.Range("I38").Value = -0.0002
$ws.Range("J38").Value = -0.0001
$ws.Range("K38").Value = -0.0001
$ws.Range("L38").Value = -0.0001
$ws.Range("M38").Value = 0
$ws.Range("I39").Value = 0.0014
$ws.Range("I40").Value = -0.0142
$ws.Range("J40").Value = -0.014
$ws.Range("K40").Value = -0.0139
$ws.Range("L40").Value = -0.0085
$ws.Range("I42").Value = 0.0027
$ws.Range("J42").Value = -0.0005
$ws.Range("K42").Value = 0.0002
$ws.Range("L42").Value = 0.0003
$ws.Range("M42").Value = 0
$ws.Range("N42").Value = 0.0001
$ws.Range("O42").Value = 0.0002
$ws.Range("P42").Value = 0.0001
$ws.Range("Q42").Value = -0.0001
$ws.Range("I43").Value = 0.002
$ws.Range("J43").Value = 0.0017
$ws.Range("K43").Value = 0.0014
$ws.Range("L43").Value = 0.0012
$ws.Range("M43").Value = 0.0006
$ws.Range("N43").Value = 0.0003
$ws.Range("O43").Value = 0.0005
$ws.Range("P43").Value = 0.0004
$ws.Range("Q43").Value = 0.0001
$ws.Range("I44").Value = 0.0217
$ws.Range("J44").Value = -0.0059
$ws.Range("K44").Value = -0.029
$ws.Range("L44").Value = -0.0072
$ws.Range("M44").Value = 0.0262
$ws.Range("N44").Value = 0.028
$ws.Range("O44").Value = 0.0309
$ws.Range("P44").Value = 0.013
$ws.Range("Q44").Value = -0.0111
$ws.Range("I46").Value = 0.0735
$ws.Range("J46").Value = -0.0728
$ws.Range("I50").Value = -0.0735
$ws.Range("J50").Value = 0.0727
$ws.Range("I51").Value = 0.0001
$ws.Range("J51").Value = 0.0001
$ws.Range("K51").Value = 0.0001
$ws.Range("L51").Value = 0.0001
$ws.Range("M51").Value = 0.0001
$ws.Range("N51").Value = 0.0001
$ws.Range("O51").Value = 0.0001
$ws.Range("P51").Value = 0.0001
$ws.Range("Q51").Value = 0.0001
$ws.Range("I52").Value = -0.0175
$ws.Range("J52").Value = -0.0179
$ws.Range("K52").Value = -0.0181
$ws.Range("L52").Value = -0.0071
$ws.Range("M52").Value = -0.003
$ws.Range("N52").Value = -0.0021
$ws.Range("O52").Value = -0.0012
$ws.Range("P52").Value = 0.0015
$ws.Range("Q52").Value = 0.0013
$ws.Range("I53").Value = 0.0389
$ws.Range("J53").Value = 0.022
$ws.Range("K53").Value = 0.022
$ws.Range("L53").Value = 0.0188
$ws.Range("M53").Value = 0.0186
$ws.Range("N53").Value = 0.0184
$ws.Range("O53").Value = 0.0183
$ws.Range("P53").Value = 0.0027
